$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new ranking snapshot row (row 19) to the tracking table
$ws.Cells.Item(19, 1).Value = "2025/12/03 11:00"
$ws.Cells.Item(19, 2).Value = "-"
$ws.Cells.Item(19, 3).Value = "-"
$ws.Cells.Item(19, 4).Value = "-"
$ws.Cells.Item(19, 5).Value = "-"
$ws.Cells.Item(19, 6).Value = "-"
$ws.Cells.Item(19, 7).Value = "-"
